$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update execution flag for "Placing COD order" row from YES to NO
$ws.Cells.Item(9, 3).Value = "NO"

# Insert two new rows before row 11 (shifts old rows 11-15 down to 13-17)
$ws.Rows.Item(11).Resize(2).Insert()

# New row 11: checking Breadcrums in listing page
$ws.Cells.Item(11, 1).Value = "checking Breadcrums in listing page"
$ws.Cells.Item(11, 2).Value = "“”"
$ws.Cells.Item(11, 3).Value = "NO"
$ws.Cells.Item(11, 4).Value = "productCatalogPage"
$ws.Cells.Item(11, 5).Value = "breadCrums"

# New row 12: checking Beadcrum in productview page
$ws.Cells.Item(12, 1).Value = "checking Beadcrum in productview page"
$ws.Cells.Item(12, 2).Value = "“”"
$ws.Cells.Item(12, 3).Value = "YES"
$ws.Cells.Item(12, 4).Value = "productCatalogPage"
$ws.Cells.Item(12, 5).Value = "productDetailPage"
$ws.Cells.Item(12, 6).Value = "breadCrums"

# Move selection to C11 to match the author's final cursor position
$ws.Range("C11").Select()

Write-Host "edit applied"
